$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The numeric-looking values in column B are stored as text in this
# workbook (no numeric formatting applied), so force column B to a text
# format before writing the new values to keep them as text, matching
# the original data type.
$ws.Range("B1:B4").NumberFormat = "@"

# Move "Starbucks" to the top of the list (row 1) with its updated value,
# shifting the remaining companies down a row; "Burger King" stays put.
$ws.Range("A1").Value = "Starbucks"
$ws.Range("B1").Value = "43600"
$ws.Range("A2").Value = "Mcdonald's"
$ws.Range("B2").Value = "37300"
$ws.Range("A3").Value = "Dell"
$ws.Range("B3").Value = "14700"
$ws.Range("A4").Value = "Burger King"
$ws.Range("B4").Value = "8360"

# Restore the default (Normal) style so the cells don't carry a leftover
# custom number format now that the text values have been committed.
$ws.Range("B1:B4").Style = "Normal"
